# Applies the "Update - Some report stuff" edits to Assignment3-report.docx
#
# wdReplaceAll = 2, wdFindContinue = 1 (wrap)

$d = $word.ActiveDocument

function Replace-Text($find, $repl) {
    $ok = $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                                   $true, 1, $false, $repl, 2)
    if (-not $ok) {
        Write-Output "WARNING: replace failed for: $find"
    }
}

# ---------------------------------------------------------------------------
# 1) "chromosome to be added" / "that chromosome is then replaced"
#    -> both become "chromosome-part"
# ---------------------------------------------------------------------------
Replace-Text "new chromosome to be added" "new chromosome-part to be added"
Replace-Text "that chromosome is then replaced" "that chromosome-part is then replaced"

# ---------------------------------------------------------------------------
# 2) "a chromosomes index simply" -> "a chromosome-parts index simply"
# ---------------------------------------------------------------------------
Replace-Text "where a chromosomes index simply" "where a chromosome-parts index simply"

# ---------------------------------------------------------------------------
# 3) "either ... or the shift" -> "neither ... nor the shift"
#    "The chromosomes FEDC" -> "The chromosome-parts FEDC"
# ---------------------------------------------------------------------------
Replace-Text "something that either the swap or the shift function" `
             "something that neither the swap nor the shift function"
Replace-Text "The chromosomes FEDC will need" "The chromosome-parts FEDC will need"

# ---------------------------------------------------------------------------
# 4) "The problem actually lies" -> "The problem might actually lie"
#    "is too good to be able" -> "is too “good” to be able"
# ---------------------------------------------------------------------------
Replace-Text "The problem actually lies in the population" `
             "The problem might actually lie in the population"

$ldq = [char]0x201C
$rdq = [char]0x201D
$findGood = "The population is too good to be able"
$replGood = "The population is too " + $ldq + "good" + $rdq + " to be able"
Replace-Text $findGood $replGood

# ---------------------------------------------------------------------------
# 5) "a chance grow freely and could improve the performance in the long
#    run, well, it's just a thought." -> "a chance to grow freely and could
#    possibly improve the solution in the long run. "
# ---------------------------------------------------------------------------
$apostrophe = [char]0x2019
$findThought = "a chance grow freely and could improve the performance in the long run, well, it" + $apostrophe + "s just a thought."
$replThought = "a chance to grow freely and could possibly improve the solution in the long run. "
Replace-Text $findThought $replThought

# ---------------------------------------------------------------------------
# 6) Resize the "Chromosome-Swap" figure (Bildobjekt 14 / InlineShapes #3)
#    from 2937600 x 1440000 EMU to 2570400 x 1260000 EMU (87.5%).
#    1 pt = 12700 EMU, so use the exact point equivalents.
# ---------------------------------------------------------------------------
$shp = $d.InlineShapes.Item(3)
$shp.Width = 2570400 / 12700
$shp.Height = 1260000 / 12700
